# Update "想去人数" (F column) figures to reflect the newly generated data
# (output generated at 456a3b4).
#
# Sheet "展览" (Worksheets.Item(1)) and sheet "全部类型" (Worksheets.Item(4))
# both list the same events, but "全部类型" has one extra row (row 8), so the
# affected rows are shifted by one starting there.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5536
$ws1.Range("F5").Value = 310
$ws1.Range("F6").Value = 831
$ws1.Range("F7").Value = 41
$ws1.Range("F8").Value = 364
$ws1.Range("F9").Value = 2
$ws1.Range("F10").Value = 3
$ws1.Range("F11").Value = 19

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5536
$ws4.Range("F5").Value = 310
$ws4.Range("F6").Value = 831
$ws4.Range("F7").Value = 41
$ws4.Range("F9").Value = 364
$ws4.Range("F10").Value = 2
$ws4.Range("F11").Value = 3
$ws4.Range("F12").Value = 19
